# Update the applicant roster: names were corrected/updated for the
# project controller's "update" export, and the Accepted/Applied/Projects
# status flags for rows 4 and 5 were swapped to reflect the corrected data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename applicants (column A, rows 2-5)
$ws.Range("A2").Value = "Audreanne Adams"
$ws.Range("A3").Value = "Will Mayer"
$ws.Range("A4").Value = "Elvie Stanton"
$ws.Range("A5").Value = "Calista Swaniawski"

# Row 4 (Elvie Stanton): Projects/Applied/Accepted -> 0,0,0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0

# Row 5 (Calista Swaniawski): Projects/Applied/Accepted -> 1,1,1
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Column A needs to widen slightly to fit the new (longer) longest name
# ("Calista Swaniawski"), matching the bestFit-recalculated width.
$ws.Columns.Item(1).ColumnWidth = 21.5
